# master-gender.xlsx: fix the Arabic "Male" translation on row 5 (B5) which
# incorrectly pointed at the Arabic string for "Female", then tidy up the
# sheet view / column sizing / page setup the way Excel would after a user
# edits that cell and glances at the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 = MLE / ara -> correct the Arabic translation to "male" (الذكر)
$ws.Range("B5").Value = "الذكر"

# User ends up with the selection sitting on D16 after the edit
$ws.Range("D16").Select()

# Column B (the translated "name" column) is sized to fit its longest entry
$ws.Columns("B:B").ColumnWidth = 6.5

# Page setup as left by Excel: Letter/A4-class paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
